$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet and update the tab/workbook title to reflect the new "through" date.
$ws.Name = "Through 2022-11-10"

# Update the row label for November to reflect the new "through" date.
$ws.Range("A12").Value = "November (through 11-10)"

# Update November row (row 12) values for columns C..I (2016..2022).
$ws.Range("C12").Value = 23
$ws.Range("D12").Value = 37
$ws.Range("E12").Value = 25
$ws.Range("F12").Value = 17
$ws.Range("G12").Value = 63
$ws.Range("H12").Value = 70
$ws.Range("I12").Value = 27

# Update Total row (row 13) values for columns C..I (2016..2022).
$ws.Range("C13").Value = 509
$ws.Range("D13").Value = 747
$ws.Range("E13").Value = 640
$ws.Range("F13").Value = 499
$ws.Range("G13").Value = 1120
$ws.Range("H13").Value = 1511
$ws.Range("I13").Value = 1426
